$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.088.29"
$ws.Range("E2").Value = "  +7.49%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.338.62"
$ws.Range("E3").Value = "  +2.49%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "412.06"
$ws.Range("E5").Value = "  +4.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "117.42"
$ws.Range("E6").Value = "  +7.63%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.333.87"
$ws.Range("E7").Value = "  +2.49%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.575"
$ws.Range("E8").Value = "  -1.78%  "
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.636"
$ws.Range("E10").Value = "  +1.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.116"
$ws.Range("E11").Value = "  +18.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "40.56"
$ws.Range("E12").Value = "  +3.04%  "
$ws.Range("E13").Value = "  -0.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.863.17"
$ws.Range("E14").Value = "  +2.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.34"
$ws.Range("E15").Value = "  -0.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.24"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.339.36"
$ws.Range("E17").Value = "  +2.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.946.58"
$ws.Range("E18").Value = "  +7.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.01"
$ws.Range("E19").Value = "  -2.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.84"
$ws.Range("E20").Value = "  +0.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.39"
$ws.Range("E21").Value = "  +1.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0000115"
$ws.Range("E22").Value = "  +5.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.49"
$ws.Range("E23").Value = "  -3.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "297.57"
$ws.Range("E24").Value = "  +0.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.22"
$ws.Range("E25").Value = "  -0.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.13"
$ws.Range("E26").Value = "  -1.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "29.24"
$ws.Range("E27").Value = "  +3.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.80"
$ws.Range("E28").Value = "  +6.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.24"
$ws.Range("E29").Value = "  -2.21%  "
$ws.Range("E30").Value = "  +1.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.51"
$ws.Range("E31").Value = "  -2.26%  "
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "43.29"
$ws.Range("E32").Value = "  +8.17%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.114"
$ws.Range("E33").Value = "  +4.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.35"
$ws.Range("E35").Value = "  +0.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.52"
$ws.Range("E36").Value = "  +18.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0494"
$ws.Range("E37").Value = "  +0.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.37"
$ws.Range("E38").Value = "  +1.83%  "
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.08"
$ws.Range("E40").Value = "  +4.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.43"
$ws.Range("E41").Value = "  -1.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "134.89"
$ws.Range("E42").Value = "  -2.13%  "
$ws.Range("E43").Value = "  +3.11%  "
$ws.Range("E44").Value = "  -0.91%  "
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.97"
$ws.Range("E46").Value = "  -1.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.44"
$ws.Range("E47").Value = "  -4.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.23"
$ws.Range("E48").Value = "  +4.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "21.21"
$ws.Range("E49").Value = "  -4.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.149.81"
$ws.Range("E50").Value = "  -0.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.662.38"
$ws.Range("E51").Value = "  +2.52%  "
